# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1470
$ws1.Range("F4").Value = 1756
$ws1.Range("F6").Value = 144
$ws1.Range("F7").Value = 656
$ws1.Range("F10").Value = 553
$ws1.Range("F16").Value = 71
$ws1.Range("F17").Value = 104
$ws1.Range("F18").Value = 4782
$ws1.Range("F20").Value = 824
$ws1.Range("F25").Value = 2069

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1470
$ws4.Range("F4").Value = 1756
$ws4.Range("F6").Value = 144
$ws4.Range("F7").Value = 656
$ws4.Range("F10").Value = 553
$ws4.Range("F16").Value = 71
$ws4.Range("F17").Value = 104
$ws4.Range("F18").Value = 4782
$ws4.Range("F22").Value = 824
$ws4.Range("F27").Value = 2069
